$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing grade ("nota") values for the last two rows of the
# "subject" table (F24 and F25). All other cells in the sheet depend on
# these two values via formulas and will recalculate automatically.
$ws.Range("F24").Value = 15
$ws.Range("F25").Value = 18

$excel.Calculate()
